$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Course_Master___12_Programs")

# New batch: "Canva for Beginners" added as row 14
$row = 14
$ws.Cells.Item($row, 1).Value = "Canva for Beginners"
$ws.Cells.Item($row, 2).Value = "CANFOR"
$ws.Cells.Item($row, 3).Value = "Graphic Design"
$ws.Cells.Item($row, 4).Value = "1 Month"
$ws.Cells.Item($row, 5).Value = 26
$ws.Cells.Item($row, 6).Value = 30
$ws.Cells.Item($row, 7).Value = 6500
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 17).Value = 5
$ws.Cells.Item($row, 18).Value = 16
$ws.Cells.Item($row, 19).Value = "Yes"
$ws.Cells.Item($row, 20).Value = "Global IT Education"
$ws.Cells.Item($row, 21).Value = "Final speech + portfolio"
$ws.Cells.Item($row, 22).Value = "Monâ€“Sat, 2 hrs/day"
$ws.Cells.Item($row, 23).Value = "Yes"
$ws.Cells.Item($row, 24).Value = "Yes"
$ws.Cells.Item($row, 25).Value = 13
$ws.Cells.Item($row, 26).Value = "Active"

# Update the view state to match the end-user's last selection/scroll
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Activate()
$ws.Range("X14").Select()
